$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(104, 1).Value = "2024-12-19 00:21:50"
$ws.Cells.Item(104, 2).Value = -0.122344933701097
$ws.Cells.Item(104, 3).Value = -0.001671089439003996
$ws.Cells.Item(104, 4).Value = 0.008177973064941894

$ws.Cells.Item(105, 1).Value = "2024-12-19 00:21:51"
$ws.Cells.Item(105, 2).Value = -0.1220716542197624
$ws.Cells.Item(105, 3).Value = -0.001715275256721996
$ws.Cells.Item(105, 4).Value = 0.008375459521211269

$ws.Cells.Item(106, 1).Value = "2024-12-19 00:21:52"
$ws.Cells.Item(106, 2).Value = -0.1227971461761626
$ws.Cells.Item(106, 3).Value = -0.001740885822317996
$ws.Cells.Item(106, 4).Value = 0.008551032431967682

$ws.Cells.Item(107, 1).Value = "2024-12-19 00:21:53"
$ws.Cells.Item(107, 2).Value = -0.1218406679914915
$ws.Cells.Item(107, 3).Value = -0.001723727755643995
$ws.Cells.Item(107, 4).Value = 0.008400805647325552

$ws.Cells.Item(108, 1).Value = "2024-12-19 00:21:54"
$ws.Cells.Item(108, 2).Value = -0.1237145844349288
$ws.Cells.Item(108, 3).Value = -0.001837709956675995
$ws.Cells.Item(108, 4).Value = 0.009094060944084072

$ws.Cells.Item(109, 1).Value = "2024-12-19 00:21:55"
$ws.Cells.Item(109, 2).Value = -0.123226585361117
$ws.Cells.Item(109, 3).Value = -0.001777530188901996
$ws.Cells.Item(109, 4).Value = 0.008761559022187769

$ws.Cells.Item(110, 1).Value = "2024-12-19 00:21:56"
$ws.Cells.Item(110, 2).Value = -0.1221367207629373
$ws.Cells.Item(110, 3).Value = -0.001912517102823994
$ws.Cells.Item(110, 4).Value = 0.009343542693678245

$ws.Cells.Item(111, 1).Value = "2024-12-19 00:21:57"
$ws.Cells.Item(111, 2).Value = -0.1235291447868803
$ws.Cells.Item(111, 3).Value = -0.001829459912817996
$ws.Cells.Item(111, 4).Value = 0.009039664738091509

$ws.Cells.Item(112, 1).Value = "2024-12-19 00:21:58"
$ws.Cells.Item(112, 2).Value = -0.1215934151274269
$ws.Cells.Item(112, 3).Value = -0.001819185318319997
$ws.Cells.Item(112, 4).Value = 0.00884803822416814

$ws.Cells.Item(113, 1).Value = "2024-12-19 00:21:59"
$ws.Cells.Item(113, 2).Value = -0.1233859983918956
$ws.Cells.Item(113, 3).Value = -0.001891816072529995
$ws.Cells.Item(113, 4).Value = 0.009336944595317927

$ws.Cells.Item(114, 1).Value = "2024-12-19 00:22:00"
$ws.Cells.Item(114, 2).Value = -0.1244042897925829
$ws.Cells.Item(114, 3).Value = -0.001952046454069995
$ws.Cells.Item(114, 4).Value = 0.0097137181104283

$ws.Cells.Item(115, 1).Value = "2024-12-19 00:22:01"
$ws.Cells.Item(115, 2).Value = -0.1233014118857682
$ws.Cells.Item(115, 3).Value = -0.001954273459773996
$ws.Cells.Item(115, 4).Value = 0.009638587072040747

$ws.Cells.Item(116, 1).Value = "2024-12-19 00:22:02"
$ws.Cells.Item(116, 2).Value = -0.122039120948175
$ws.Cells.Item(116, 3).Value = -0.001917780934487995
$ws.Cells.Item(116, 4).Value = 0.009361771976643377

$ws.Cells.Item(117, 1).Value = "2024-12-19 00:22:03"
$ws.Cells.Item(117, 2).Value = -0.1229988524600049
$ws.Cells.Item(117, 3).Value = -0.001889892749421995
$ws.Cells.Item(117, 4).Value = 0.009298185578055557

$ws.Cells.Item(118, 1).Value = "2024-12-19 00:22:04"
$ws.Cells.Item(118, 2).Value = -0.1241440236198832
$ws.Cells.Item(118, 3).Value = -0.001860334310077996
$ws.Cells.Item(118, 4).Value = 0.009237975461248075

$ws.Cells.Item(119, 1).Value = "2024-12-19 00:22:05"
$ws.Cells.Item(119, 2).Value = -0.1218666946087615
$ws.Cells.Item(119, 3).Value = -0.001825410811537996
$ws.Cells.Item(119, 4).Value = 0.008898271276209298

$ws.Cells.Item(120, 1).Value = "2024-12-19 00:22:06"
$ws.Cells.Item(120, 2).Value = -0.1215608818558394
$ws.Cells.Item(120, 3).Value = -0.001858208531905997
$ws.Cells.Item(120, 4).Value = 0.009035418712421505

$ws.Cells.Item(121, 1).Value = "2024-12-19 00:22:07"
$ws.Cells.Item(121, 2).Value = -0.1242351167803281
$ws.Cells.Item(121, 3).Value = -0.001881642705563995
$ws.Cells.Item(121, 4).Value = 0.009350644050583818

$ws.Cells.Item(122, 1).Value = "2024-12-19 00:22:08"
$ws.Cells.Item(122, 2).Value = -0.1230671723303385
$ws.Cells.Item(122, 3).Value = -0.001855627229839996
$ws.Cells.Item(122, 4).Value = 0.009134671843023497

$ws.Cells.Item(123, 1).Value = "2024-12-19 00:22:09"
$ws.Cells.Item(123, 2).Value = -0.1226312264910666
$ws.Cells.Item(123, 3).Value = -0.001787906010931995
$ws.Cells.Item(123, 4).Value = 0.008770124278853634

$ws.Cells.Item(124, 1).Value = "2024-12-19 00:22:10"
$ws.Cells.Item(124, 2).Value = -0.1226051998737967
$ws.Cells.Item(124, 3).Value = -0.001771456536981996
$ws.Cells.Item(124, 4).Value = 0.008687591311376853

$ws.Cells.Item(125, 1).Value = "2024-12-19 00:22:11"
$ws.Cells.Item(125, 2).Value = -0.122702799688559
$ws.Cells.Item(125, 3).Value = -0.001768267869723995
$ws.Cells.Item(125, 4).Value = 0.008678856728578335

$ws.Cells.Item(126, 1).Value = "2024-12-19 00:22:12"
$ws.Cells.Item(126, 2).Value = -0.122549893312098
$ws.Cells.Item(126, 3).Value = -0.001775201955665996
$ws.Cells.Item(126, 4).Value = 0.008702032410971819
